# B1--and-B2-PowerPoint.pptx edit
#  1) Table on slide 5 (graphicFrame #2) gets a new built-in table style.
#  2) The presentation's theme colour scheme switches from the
#     "Integral" (Red Violet) palette to the default "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Apply the new table style to the table on slide 5 -------------
$slide5 = $p.Slides.Item(5)
$tableShape = $slide5.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{6BB579F5-9EE0-4D3E-AB9E-9028A1B785BD}")

# --- 2. Switch the theme colour scheme to the Office palette ----------
$colorScheme = $slide5.ThemeColorScheme

$colorScheme.Item(1).RGB  = 0         # dk1      000000
$colorScheme.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 6968388   # dk2      44546A
$colorScheme.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colorScheme.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colorScheme.Item(6).RGB  = 3243501   # accent2  ED7D31
$colorScheme.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colorScheme.Item(8).RGB  = 49407     # accent4  FFC000
$colorScheme.Item(9).RGB  = 12874308  # accent5  4472C4
$colorScheme.Item(10).RGB = 4697456   # accent6  70AD47
$colorScheme.Item(11).RGB = 12673797  # hlink    0563C1
$colorScheme.Item(12).RGB = 7491477   # folHlink 954F72
